$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# BGT and WM voxel processing in a loop:
# A new voxel record (LWP758) is appended as row 23, mirroring the
# structure of the existing rows. The LWP run label in the BA1 header
# cell is also bumped to match the newest voxel batch.

$ws.Range("BA1").Value = "LWP758"

$newRow = 23
$values = @(
    "LWP758",
    "3",
    "5",
    "1",
    "2.176556e+00",
    "0.000000e+00",
    "7.614723e+00",
    "1.061412e+01",
    "0.000000e+00",
    "0.000000e+00",
    "5.222804e+00",
    "1.385666e+01",
    "9.927932e+00",
    "4.419998e+00",
    "0.000000e+00",
    "1.749804e+00",
    "1.089423e+01",
    "4.949827e+00",
    "4.461056e+00",
    "8.102939e-01",
    "0.000000e+00",
    "1.007113e+01",
    "3.971859e+00",
    "1.584405e+01",
    "4.461056e+00",
    "2.378460e+01",
    "LWP758",
    "3",
    "5",
    "1",
    "3.243683e+00",
    "1.661905e+00",
    "7.694937e+00",
    "2.774940e+00",
    "3.303601e+00",
    "2.200221e+00",
    "8.348805e+00",
    "9.336860e+00",
    "6.189578e+00",
    "3.951900e+00",
    "6.515631e+00",
    "3.200276e+00",
    "3.863494e+00",
    "5.240933e+00",
    "2.324456e+00",
    "1.662289e+00",
    "5.343794e-01",
    "5.839715e+00",
    "2.880158e+00",
    "3.248666e+00",
    "1.097800e+00",
    "9.803526e+00",
    "LWP758",
    "3",
    "5",
    "1",
    "1.332186e+00",
    "8.427274e+01",
    "9.760060e+01",
    "0.000000e+00",
    "4.393711e-02",
    "5.615234e+00",
    "2.146096e+00",
    "2.858998e+00",
    "2.858998e+00",
    "1.119265e-02",
    "1.219337e-02",
    "4.629942e+00",
    "1.000000e+00",
    "3.585189e+01",
    "2.483205e-02",
    "-3.859145e+01",
    "-2.000000e+01",
    "2.087835e+01",
    "1.062929e+01",
    "8.317021e-02",
    "2.075195e+00",
    "1.712270e+03",
    "2.038166e-01",
    "6.713867e+00",
    "4.637078e-02",
    "9.564863e-02",
    "-1.116521e-01",
    "-4.917817e-03",
    "7.367890e+00",
    "6.932134e+00"
)

# Force the whole new row to be stored as text (matches the text-typed
# cells used throughout the sheet for this voxel data) before writing
# the values in a loop across columns A:CH.
$rowRange = $ws.Range($ws.Cells.Item($newRow, 1), $ws.Cells.Item($newRow, $values.Length))
$rowRange.NumberFormat = "@"

for ($col = 1; $col -le $values.Length; $col++) {
    $ws.Cells.Item($newRow, $col).Value = $values[$col - 1]
}
